# Work in progress: Created a parametrized python script which contains all
# the methods of running the fuzzer, eliminating the need of having three
# redundant scripts (ScriptFullPathOfflineAll / ScriptFullPathDFSSymbolic /
# ScriptFullPathConcolic) -> single "ScriptFullPath".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 6 (ScriptWorkingFolder and its value) is unchanged.

# Row 10: RobotModelFullPath moves up to row 10 (from row 14), now with its description
$ws.Range("A10").Value = "RobotModelFullPath"
$ws.Range("B10").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4\Applications\C#Models\SimpleBankLoanCSharp"
$ws.Range("C10").Value = "*Full path of the robot model under test"

# Row 8: consolidate ScriptFullPathOfflineAll into the new single ScriptFullPath
$ws.Range("A8").Value = "ScriptFullPath"
$ws.Range("B8").Value = "C:\GithubPhD\rpa-testing\TestingTool_v4"

# Row 4: PyhtonHomePath value changed (different machine/user)
$ws.Range("B4").Value = "C:\Users\Marina Cernat\AppData\Local\Programs\Python\Python38"

# Remove the now-redundant rows that used to hold
# ScriptFullPathDFSSymbolic / ScriptFullPathConcolic / old RobotModelFullPath
$ws.Range("A11:C14").Clear()

# Remove the trailing blank rows that are no longer present
$ws.Rows("1002:1005").Delete()

# Update selection on the Settings sheet
$ws.Range("B8").Select()
